# #5: property boat&car done
# Fix the "汽車" (car) sheet: row 1 was accidentally populated with a
# duplicate of the first data row instead of real column headers, and the
# data rows were missing the trailing property_category/category/date/
# legislator_name/legislator_id/source_file/index columns that every other
# sheet in this workbook already has. Rebuild row 1 as headers and fill in
# the missing columns (H:N) for the two existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Row 1: real column headers (was a stray copy of row 2's data) ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the header styling (bold, centered, thin border) already used on
# B1:G1 so the newly-added H1:N1 header cells look the same.
$headerRange = $ws.Range("H1:N1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Row 2 (car #46 - NISSAN CEFIRO) ----
# Leading apostrophe forces the "2012-04-19" literal to stay text instead of
# being auto-converted to a date serial by Excel's input parser.
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "'2012-04-19"
$ws.Range("K2").Value = "李慶華"
$ws.Range("L2").Value = 607
$ws.Range("M2").Value = "tmpe2cb1"
$ws.Range("N2").Value = 46

# ---- Row 3 (car #47 - 中華 GRUNDER) ----
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "'2012-04-19"
$ws.Range("K3").Value = "李慶華"
$ws.Range("L3").Value = 607
$ws.Range("M3").Value = "tmpe2cb1"
$ws.Range("N3").Value = 47
